$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Table style swap on slide 6 (the sources-of-finance table):
#    {F182BE06-0029-4D08-B2A1-942E09A6F6A6} -> {35387266-C5AC-41E1-B970-BFE7E2CAFED0}
# ---------------------------------------------------------------------------
$slide6 = $p.Slides.Item(6)
$tableShape = $slide6.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{35387266-C5AC-41E1-B970-BFE7E2CAFED0}")

# ---------------------------------------------------------------------------
# 2) Theme colour swap: the deck's main theme ("Integral") is replaced with
#    the stock "Office Theme" palette. dk1/lt1 (black/white) are unchanged;
#    the rest of the 12 theme colours move to the standard Office values.
# ---------------------------------------------------------------------------
function Set-ThemeColor($themeColors, [int]$index, [string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $rgb = $r + ($g * 256) + ($b * 65536)
    $color = $themeColors.Colors($index)
    $color.RGB = $rgb
}

$themeColors = $p.Slides.Item(1).ThemeColorScheme

Set-ThemeColor $themeColors 1  "000000"   # Dark 1
Set-ThemeColor $themeColors 2  "FFFFFF"   # Light 1
Set-ThemeColor $themeColors 3  "44546A"   # Dark 2
Set-ThemeColor $themeColors 4  "E7E6E6"   # Light 2
Set-ThemeColor $themeColors 5  "5B9BD5"   # Accent 1
Set-ThemeColor $themeColors 6  "ED7D31"   # Accent 2
Set-ThemeColor $themeColors 7  "A5A5A5"   # Accent 3
Set-ThemeColor $themeColors 8  "FFC000"   # Accent 4
Set-ThemeColor $themeColors 9  "4472C4"   # Accent 5
Set-ThemeColor $themeColors 10 "70AD47"   # Accent 6
Set-ThemeColor $themeColors 11 "0563C1"   # Hyperlink
Set-ThemeColor $themeColors 12 "954F72"   # Followed Hyperlink
